# ============================================================================
# Applies the "Queries" worksheet addition to AgentInteractionReportData.xlsx
# (commit: "commited on april 15th")
# ============================================================================

$wb = $excel.ActiveWorkbook

# --- Minor selection changes left behind on pre-existing sheets ---
$sheetShow = $wb.Worksheets.Item("Show")
$sheetShow.Activate()
$sheetShow.Range("B14").Select()

$sheetAdvanceSearch = $wb.Worksheets.Item("AdvanceSearch")
$sheetAdvanceSearch.Activate()
$sheetAdvanceSearch.Range("D14").Select()

# --- Add the new "Queries" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Queries"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Report Channel"
$ws.Range("B1").Value = "Report Name"
$ws.Range("C1").Value = "Report Type"
$ws.Range("D1").Value = "Start Date"
$ws.Range("E1").Value = "End Date"

# --- Data row (row 2) labels ---
$ws.Range("A2").Value = "Agent"
$ws.Range("B2").Value = "OCM Agent Interaction Report"
$ws.Range("C2").Value = "Date Range"

# Header cells A1:E1 and A2:B2 use the existing "text" style (numFmtId 49)
$ws.Range("A1:E1").NumberFormat = "@"
$ws.Range("A2:B2").NumberFormat = "@"

# New "Query" column header (insertion order matters for shared-string indices)
$ws.Range("F1").Value = "Query"

# First date-like text value -- force a leading quote-prefix so Excel keeps it as text
$ws.Range("D2").Formula = "'15-04-2020 00:00:00"

# Large SQL query text, wrapped and centered
$ws.Range("F2").Value = 'SELECT M.[AgentID] as [Agent ID],A.[AgentName] as [Agent Name],A.[SupervisorName] as [Supervisor Name],[Channel],[Direction],[DNIS] as [Local Party],' + "`r`n" + '[Ani] as [Remote Party],FORMAT([dbo].[VARCHARTODATETIME](M.ConnectedDateTime),''dd/MM/yyyy HH:mm:ss'') as [Interaction Connected Date Time],' + "`r`n" + 'FORMAT([dbo].[VARCHARTODATETIME](M.[CreatedDateTime]),''dd/MM/yyyy HH:mm:ss'') as [Created Date Time],[SessionID] as [Session ID],CONVERT(varchar, DATEADD(ms, M.QueueTime* 1000, 0), 108) as [Queue Time],' + "`r`n" + 'CONVERT(varchar, DATEADD(ms, M.ActiveTime* 1000, 0), 108) as [Active Time],' + "`r`n" + 'CONVERT(varchar, DATEADD(ms, M.HoldTime* 1000, 0), 108) as [Hold Time],' + "`r`n" + 'CONVERT(varchar, DATEADD(ms, M.ACWTime* 1000, 0), 108) as [ACW Time],' + "`r`n" + 'CONVERT(varchar, DATEADD(ms, M.HandleTime* 1000, 0), 108) as [Handle Time],' + "`r`n" + 'case when [IsConferenced]=1 Then ''true'' when [IsConferenced]=0 then ''false'' end as [Is Conferenced],' + "`r`n" + 'case when [IsTransfered]=1 Then ''true'' when [IsTransfered]=0 then ''false'' end as [Is Transfered], ' + "`r`n" + '[TPINTransferReconnected] as [TPIN Transfer Reconnected],[SubChannel] as [Sub Channel],' + "`r`n" + '[SubSessionID] as [Sub Session ID],[InteractionID] as [Interaction ID],[Skill],[SkillName] as [Skill Name],' + "`r`n" + '[DNISName] as [DNIS Name],[TransferedTo] as [Transfered To],[ConferencedTo] as [Conferenced To],' + "`r`n" + '[ConferenceToAgentList] as [Conference To Agent List],[TransferToAgent] as [Transfer To Agent],' + "`r`n" + '[TransferConferenceFromAgent] as [Transfer Conference From Agent] ,[TransferConferenceFromInteraction] as [Transfer Conference From Interaction],' + "`r`n" + '--[OtherData],' + "`r`n" + 'FORMAT([dbo].[VARCHARTODATETIME](M.[ClosedDateTime]),''dd/MM/yyyy HH:mm:ss'') as [Closed Date Time],' + "`r`n" + 'FORMAT([dbo].[VARCHARTODATETIME](M.[DisconnectedDateTime]),''dd/MM/yyyy HH:mm:ss'') as [Interaction Disconnected Date Time],[ClosedReason] as [Closed Reason],' + "`r`n" + '[CIF],[RegisteredMobileNo] as [Registered Mobile No] FROM (SELECT DISTINCT [User]  AS Ani,' + "`r`n" + '							AgentId as AgentID,' + "`r`n" + '							T.Channel,' + "`r`n" + '							SubChannel,' + "`r`n" + '							T.SessionID AS SessionID,' + "`r`n" + '							SubSessionId as SubSessionID,' + "`r`n" + '							InteractionId as InteractionID,' + "`r`n" + '							T.Direction,' + "`r`n" + '							CreatedDateTime,' + "`r`n" + '							CreatedReason,' + "`r`n" + '							Skill,' + "`r`n" + '							TS.SkillName,' + "`r`n" + '							ISNULL(A.FirstName,'''') +'' ''+ ISNULL(A.LastName,'''') AS AgentName,' + "`r`n" + '							Dnis as DNIS,' + "`r`n" + '							DnisName as DNISName,' + "`r`n" + '							IsTransfered,' + "`r`n" + '							IsConferenced,' + "`r`n" + '							IsReconnected AS TPinTransferReconnected,' + "`r`n" + '							IsConferencedTo AS ConferencedTo,' + "`r`n" + '							IsTranferedTo AS TransferedTo,' + "`r`n" + '							CASE WHEN IsTransfered=1 OR IsConferenced=1 THEN  TrasnferConferenceFromAgent ELSE '''' END AS TransferConferenceFromAgent,' + "`r`n" + '							CASE WHEN IsTransfered=1 OR IsConferenced=1 THEN  TrasnferConferenceFromInteraction ELSE '''' END AS TransferConferenceFromInteraction,' + "`r`n" + '							OtherData,' + "`r`n" + '							ClosedDateTime AS ClosedDateTime,' + "`r`n" + '							ClosedReason,' + "`r`n" + '							CallConnectedTime AS  ConnectedDateTime,' + "`r`n" + '							CallDisconnectedTime AS  DisconnectedDateTime,' + "`r`n" + '							ActiveTime,' + "`r`n" + '							HoldTime,' + "`r`n" + '							TrasnferToAgent AS TransferToAgent,' + "`r`n" + '							ConferenceToAgentList,' + "`r`n" + '							QueueTime,' + "`r`n" + '							AcwTime as ACWTime,' + "`r`n" + '							ActiveTime+HoldTime+AcwTime HandleTime,' + "`r`n" + '							IH.CIF,' + "`r`n" + '							IH.CLID AS RegisteredMobileNo' + "`r`n" + '							FROM TMAC_Interactions T with(nolock)' + "`r`n" + '							INNER JOIN AGT_Agent A with(nolock) ON A.AvayaLoginID=T.AgentId ' + "`r`n" + '							LEFT JOIN AGT_Agent AA with(nolock) ON AA.AvayaLoginID = T.TrasnferConferenceFromAgent' + "`r`n" + '							LEFT JOIN GBL_InteractionHistory IH WITH(NOLOCK) ON IH.SessionID=T.SessionId ' + "`r`n" + '							LEFT JOIN TMAC_Skills TS WITH(NOLOCK) ON TS.SkillExtension=T.Skill' + "`r`n" + '							where 1=1 AND IH.ID IN (SELECT MIN(ID) FROM GBL_InteractionHistory WHERE SESSIONID=T.SessionId AND CLID IS NOT NULL)' + "`r`n" + '							AND  ClosedDateTime>=''ReportBeforeDate'' AND ClosedDateTime<=''ReportAfterDate'') M' + "`r`n" + '                            INNER JOIN  fn_AgentHierarchy(''na'',''1'',''1'') A  ON A.AgentId=M.AgentID' + "`r`n" + '							order by M.[ClosedDateTime]'
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4160
$ws.Range("F2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.5

# Second date-like text value
$ws.Range("E2").Formula = "'15-04-2020 15:30:00"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 13.83
$ws.Columns.Item(2).ColumnWidth = 26.33
$ws.Columns.Item(3).ColumnWidth = 10.83
$ws.Columns.Item(4).ColumnWidth = 18.33
$ws.Columns.Item(5).ColumnWidth = 17.83
$ws.Columns.Item(6).ColumnWidth = 20.83

# --- Final selection / activation: Queries is the active sheet ---
$ws.Activate()
$ws.Range("E2").Select()
